$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B (pushes old column B -> column C).
# Excel's default insert picks up formatting from the column to the left (A),
# matching the style pattern seen in the target file (B gets same style as A).
$ws.Columns("B:B").Insert()

# Match column A's width on the newly inserted column B
$ws.Columns("B:B").ColumnWidth = $ws.Columns("A:A").ColumnWidth

# Update header row
$ws.Range("A1").Value = "ratingName_EN"
$ws.Range("B1").Value = "ratingName_CN"

# Update the new Chinese-language column with localized picture paths
$ws.Range("B2").Value = "Instructions_CN/ratingCS+1.png"
$ws.Range("B3").Value = "Instructions_CN/ratingCS+3.png"
$ws.Range("B4").Value = "Instructions_CN/ratingCS+4.png"

# Match final UI selection state
$ws.Range("C10").Select()
